$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.748.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.333.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.88%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -3.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.14"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.29"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.04%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.21%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.26"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0920"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.55"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.27%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.01"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.92%  "

$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.108"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.52"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.687.75"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.330.58"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.730.96"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.57"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.76%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.98"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.36"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.08%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.20"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.31%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.22%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.65"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +10.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.14"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.68%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.87"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +9.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.56"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "168.11"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0902"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.21%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +8.43%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.87%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.116"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.71"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0365"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.64%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.91"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +10.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.79"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.78%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +9.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.60"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +12.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.71"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +14.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.240"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "71.87"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.38%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "115.00"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.219"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +17.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.658.29"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.99"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.58%  "

$ws.Range("B51").Value = "MinaProtocolToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.55"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +9.81%  "

Write-Host "Applied cryptos update"
